$wb = $excel.ActiveWorkbook

# --- Step 1: insert a new sheet "2022-Q1" right before the "总计" sheet ---
# NOTE: capturing "Worksheets.Item(2)" and then calling Add(before:=that handle)
# leaves the captured variable pointing at the newly-inserted sheet (the handle
# tracks tab position, not sheet identity) - so the "总计" sheet must be
# re-acquired by its new position AFTER the insert, not reused from before it.
$beforeTarget = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeTarget)
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item(3)

# Pick up the same header / row-index cell styling already used on the
# "总计" sheet (bold, centered, thin-bordered) by copying a formatted cell
# onto the new ranges before filling in values.
$totalSheet.Range("B1").Copy($newSheet.Range("B1:H1"))
$totalSheet.Range("A2").Copy($newSheet.Range("A2:A4"))

# Header row (row 1)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Data rows (rows 2-4): A=index(number), B=code(text), C=name(text),
# D=scale(text), E=position(text), F=ratio(text), G=value(text), H=rank(number)
$data = @(
    @(0, "501030", "汇添富中证环境治理指数（LOF）A", "6.61", "93.20", "2.11", "0.1395", 6),
    @(1, "501031", "汇添富中证环境治理指数（LOF）C", "2.74", "93.20", "2.11", "0.0578", 6),
    @(2, "164908", "交银施罗德中证环境治理指数（LOF）", "2.12", "93.72", "2.16", "0.0458", 6)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $data[$r]

    $newSheet.Cells.Item($rowNum, 1).Value = $rowData[0]

    # Text columns that look numeric must be forced to Text so Excel doesn't
    # silently reinterpret "93.20" / "501030" (etc.) as a number and lose the
    # literal formatting / leading zeros.
    foreach ($pair in @(@(2,1), @(3,2), @(4,3), @(5,4), @(6,5), @(7,6))) {
        $col = $pair[0]
        $idx = $pair[1]
        $cell = $newSheet.Cells.Item($rowNum, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$idx]
    }

    $newSheet.Cells.Item($rowNum, 8).Value = $rowData[7]
}

# --- Step 2: update the "总计" sheet (now at position 3) with a new 2022-Q1 summary row ---
$totalSheet.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the header row above for
# B2:D2 (plain data cells should carry no special style) and drops the
# row-index style that used to live on A2 entirely - tidy both up: clear the
# plain cells back to the default style, and re-copy the index-column style
# from A3 (which kept it after the shift) onto the new A2.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.24

# fix the row-index column for the row that shifted down (was 2021-Q3 / index 0)
$totalSheet.Cells.Item(3, 1).Value = 1
